$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a "Total" row at row 4, directly below the single
# data row (row 3). Five new transaction rows need to be inserted above the
# Total row (which then becomes row 9), and the Total row's SUM formulas
# need to grow to cover the newly added rows.

# Insert 5 blank rows at row 4, pushing the existing "Total" row down to row 9.
$ws.Range("A4:A8").EntireRow.Insert()

# The new rows should carry the same formatting as row 3 (the existing data
# row): copy its style down into rows 4-8.
$ws.Range("A3:I3").Copy()
$ws.Range("A4:I8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 4: Funds entry
$ws.Range("A4").Value = "31-8-2022"
$ws.Range("E4").Value = "Angel"
$ws.Range("F4").Value = 2000
$ws.Range("G4").Value = "Gpay"

# Row 5: Funds entry
$ws.Range("A5").Value = "31-8-2022"
$ws.Range("E5").Value = "Angel"
$ws.Range("F5").Value = 1700
$ws.Range("G5").Value = "Gpay"

# Row 6: Expense entry
$ws.Range("A6").Value = "31-8-2022"
$ws.Range("B6").Value = "Food"
$ws.Range("C6").Value = 25
$ws.Range("D6").Value = "Cash"

# Row 7: Funds entry
$ws.Range("A7").Value = "31-8-2022"
$ws.Range("E7").Value = "Angel"
$ws.Range("F7").Value = 400
$ws.Range("G7").Value = "Card"

# Row 8: Income entry
$ws.Range("A8").Value = "31-8-2022"
$ws.Range("H8").Value = "Salary"
$ws.Range("I8").Value = 15000

# Row 9 (the "Total" row, originally row 4): extend the SUM ranges to cover
# the newly added data rows 3-8.
$ws.Range("C9").Formula = "=SUM(C3:C8)"
$ws.Range("F9").Formula = "=SUM(F3:F8)"
$ws.Range("I9").Formula = "=SUM(I3:I8)"
